# Rewrite the "Micro results" answer cell: replace the single long
# bold(=0) paragraph (many lines joined with manual line breaks) with a
# leading blank paragraph followed by one paragraph per line, each
# styled Times New Roman / 20 half-points / blue (0000FF), no bold.

$newLines = @(
    '23/09 – TB CULTURE – TISSUE;LungLowrLob R MYCOBACTERIAL CULTURE: NO GROWTH AFTER 8 WEEKS',
    '23/09 – REF FUNGAL PCR – **No clear Result** +',
    '  Summary: No findings reported.',
    '23/09 – FUNGUS CULTURE – **Negative**',
    '  Summary: Fungal culture negative.',
    '23/09 – ANAEROBIC CULTURE – **Negative**',
    '  Summary: No anaerobes isolated.',
    '23/09 – NOCARDIA CULTURE – **Negative**',
    '  Summary: No Nocardia growth.',
    '23/09 – RESP. CULT AND MICRO – TISSUE;LungLowrLob R NO GROWTH AFTER 6 DAYS',
    '23/09 – GRAM SMEAR – **Negative**',
    '  Summary: No organisms seen.',
    '23/09 – TB MICROSCOPY – **Negative**',
    '  Summary: No mycobacteria seen.',
    '22/09 – BLC – ;Femoral R NO GROWTH AFTER 5 DAYS',
    '22/09 – BLC – ;Other (specify site in Clinical Details) NO GROWTH AFTER 5 DAYS',
    '20/09 – BLC – PERIPHERAL–RIGHT NO GROWTH AFTER 5 DAYS',
    '20/09 – BLC – ;Other (specify site in Clinical Details) NO GROWTH AFTER 5 DAYS',
    '--------Previous result (1 year)--------',
    '24/07 – EBV VCA IgG – Positive',
    '19/06 – RESPIRATORY PCR – **Positive**  ',
    '**Summary:** Human Rhinovirus/Enterovirus detected.',
    '12/06 – EBV VCA IgG – Positive',
    '20/02 – RESPIRATORY PCR – **Positive**  ',
    '**Summary:** RSV detected by PCR.',
    '04/02 – RESPIRATORY PCR – **Positive**  ',
    '**Summary:** RSV detected in throat swab.',
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "Micro results" row, and a "clean" donor row (one whose
# answer cell already has no explicit <w:pPr>, e.g. "Vital signs
# rules") by scanning column 1 labels -- more robust than hard-coding
# row indices.
$targetRow = 0
$cleanRow = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $label = $t.Cell($r, 1).Range.Text
    if ($label -like "*Micro results*") {
        $targetRow = $r
    }
    if ($label -like "*Vital signs rules*") {
        $cleanRow = $r
    }
}

if ($targetRow -eq 0) {
    throw "Could not find the Micro results row"
}
if ($cleanRow -eq 0) {
    throw "Could not find the Vital signs rules row"
}

$cell = $t.Cell($targetRow, 2)

# The existing answer paragraph carries an explicit <w:pPr><w:jc val="left"/></w:pPr>.
# Copy in the paragraph format from a cell that has no explicit jc (the
# "Vital signs rules" answer) so the rewritten paragraphs come out with
# no <w:pPr> at all, matching the target markup.
$cleanFormat = $t.Cell($cleanRow, 2).Range.ParagraphFormat
$cell.Range.ParagraphFormat = $cleanFormat

# Build the full replacement text: each data line, separated by manual
# paragraph marks (CR).
$cell = $t.Cell($targetRow, 2)
$cell.Range.Text = ($newLines -join ([string][char]13))

# Insert one more blank leading paragraph ahead of everything -- this
# yields a plain <w:p><w:r/></w:p> (an empty, unformatted run) exactly
# like the sibling "rule-based" answer cells in this table.
$cell = $t.Cell($targetRow, 2)
$cell.Range.InsertParagraphBefore()

# Re-acquire the cell and compute the new cell start so the offsets
# below are correct after the text swap + paragraph insert.
$cell = $t.Cell($targetRow, 2)
$cellStart = $cell.Range.Start

# Walk the data paragraphs (skip the blank leading one) and set their
# run formatting while excluding the trailing paragraph-mark character,
# so Word doesn't stamp a <w:pPr><w:rPr> "mark" style on them.
$offset = $cellStart + 1
for ($i = 0; $i -lt $newLines.Count; $i++) {
    $line = $newLines[$i]
    $len = $line.Length
    $r = $d.Range($offset, $offset + $len)
    $r.Font.Name = "Times New Roman"
    $r.Font.NameAscii = "Times New Roman"
    $r.Font.Size = 10
    $r.Font.Bold = $false
    $r.Font.Color = 16711680
    $offset = $offset + $len + 1
}

Write-Output "Micro results cell rewritten"

